$d = $word.ActiveDocument
$sel = $word.Selection

# The cursor/selection starts at the document's last-edit location (the
# existing "_GoBack" bookmark), right after the two "123" runs. Split
# that spot into two more paragraphs, as if the author pressed Enter and
# kept typing: "123<enter>123<enter>1234".
$sel.TypeParagraph()
$sel.TypeText("123")
$sel.TypeParagraph()
$sel.TypeText("1234")

# Word also re-homes the hidden "_GoBack" bookmark to the most recent
# edit location. Recreate it collapsed right after the "1234" we just
# typed: bookmark a trailing placeholder character, then delete that
# character so the bookmark collapses in place.
$sel.TypeText("X")
$endPos = $d.Content.End - 2
$tmpRange = $d.Range($endPos, $endPos + 1)

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $tmpRange)

$d.Range($endPos, $endPos + 1).Delete()
